# Data for xG trendline
# The "Season" column (F) previously held numeric season indices (1, 2, 3).
# Replace them with the actual season labels as text, matching each block
# of rows to the season it represents:
#   rows 2-39   (old value 1) -> "19/20"
#   rows 40-77  (old value 2) -> "20/21"
#   rows 78-115 (old value 3) -> "21/22"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2:F39").Value = "19/20"
$ws.Range("F40:F77").Value = "20/21"
$ws.Range("F78:F115").Value = "21/22"

# Reflect the author's final cursor/selection position in the sheet view.
$ws.Activate()
$ws.Range("M40").Select()
